$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 05:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 587155
$ws.Range("C4").Value = 214
$ws.Range("E4").Value = 526563
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 23644

# Row 27 - Japon
$ws.Range("B27").Value = 7645
$ws.Range("C27").Value = 27
$ws.Range("E27").Value = 6703

# Row 33 - Australia
$ws.Range("B33").Value = 6366
$ws.Range("C33").Value = 7
$ws.Range("E33").Value = 2811

# Row 36 - Pakistan
$ws.Range("B36").Value = 5707
$ws.Range("C36").Value = 211
$ws.Range("D36").Value = 1097
$ws.Range("E36").Value = 4514
$ws.Range("G36").Value = 3
$ws.Range("H36").Value = 96

# Row 114 - Vietnam
$ws.Range("D114").Value = 155
$ws.Range("E114").Value = 110
